$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1500
$ws.Range("C2").Value = 5500

$ws.Range("B3").Value = 2000
$ws.Range("C3").Value = 1500

$ws.Range("B4").Value = 1500
$ws.Range("C4").Value = 5500

$ws.Range("B5").Value = 2000
$ws.Range("C5").Value = 1500

$ws.Range("B6").Value = 1800
$ws.Range("C6").Value = 5500

$ws.Range("B7").Value = 1800
$ws.Range("C7").Value = 1400

$ws.Range("B8").Value = 1800
$ws.Range("C8").Value = 1400
